# "update with new column" - Internship_list.xlsx
#
# Row 6 (Amazon) had its NOTES/link column (D6) showing a stale, copy-pasted
# URL ("http://microsoft.com/" - the same text as the Microsoft row above
# it). Fix it to point at Amazon's own internship listing page.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "https://www.amazon.jobs/en/teams/internships-for-students"

# The longer text reflows the wrapped NOTES column, and the surrounding
# rows' autofit heights were re-measured on save as well.
$ws.Rows(5).RowHeight = 34
$ws.Rows(6).RowHeight = 170
$ws.Rows(7).RowHeight = 204

# Selection/viewport as left at save time.
$ws.Range("F6").Select()
